$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price/Volume columns keep their text formatting (values like "67.470.03"
# are not valid numbers and must stay as text, matching the original inlineStr cells).
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '67.677.82'
$ws.Range("E2").Value = '  +1.93%  '

$ws.Range("D3").Value = '3.602.10'
$ws.Range("E3").Value = '  +0.85%  '

$ws.Range("E4").Value = '  -0.05%  '

$ws.Range("D5").Value = '202.01'
$ws.Range("E5").Value = '  +8.37%  '

$ws.Range("D6").Value = '560.13'
$ws.Range("E6").Value = '  -4.94%  '

$ws.Range("D7").Value = '3.597.02'
$ws.Range("E7").Value = '  +0.85%  '

$ws.Range("D8").Value = '0.614'
$ws.Range("E8").Value = '  -0.07%  '

$ws.Range("E9").Value = '  +0.21%  '

$ws.Range("D10").Value = '0.671'
$ws.Range("E10").Value = '  +0.06%  '

$ws.Range("D11").Value = '59.93'
$ws.Range("E11").Value = '  +11.74%  '

$ws.Range("D12").Value = '0.152'
$ws.Range("E12").Value = '  +3.50%  '

$ws.Range("D13").Value = '0.0000287'
$ws.Range("E13").Value = '  +11.10%  '

$ws.Range("D14").Value = '9.98'
$ws.Range("E14").Value = '  +2.20%  '

$ws.Range("D15").Value = '4.174.86'
$ws.Range("E15").Value = '  +0.84%  '

$ws.Range("D16").Value = '3.592.22'
$ws.Range("E16").Value = '  +0.58%  '

$ws.Range("E17").Value = '  +0.61%  '

$ws.Range("D18").Value = '19.01'
$ws.Range("E18").Value = '  +3.98%  '

$ws.Range("D19").Value = '67.515.35'
$ws.Range("E19").Value = '  +1.89%  '

$ws.Range("D20").Value = '12.31'
$ws.Range("E20").Value = '  +0.67%  '

$ws.Range("D21").Value = '1.08'
$ws.Range("E21").Value = '  +1.81%  '

$ws.Range("D22").Value = '402.68'
$ws.Range("E22").Value = '  +1.42%  '

$ws.Range("D23").Value = '12.76'
$ws.Range("E23").Value = '  +13.75%  '

$ws.Range("D24").Value = '4.14'
$ws.Range("E24").Value = '  -5.54%  '

$ws.Range("D25").Value = '85.21'
$ws.Range("E25").Value = '  -0.47%  '

$ws.Range("D26").Value = '2.94'
$ws.Range("E26").Value = '  +2.04%  '

$ws.Range("D27").Value = '12.51'
$ws.Range("E27").Value = '  +0.19%  '

$ws.Range("D28").Value = '3.88'
$ws.Range("E28").Value = '  +9.18%  '

$ws.Range("E29").Value = '  +1.07%  '

$ws.Range("D30").Value = '8.35'
$ws.Range("E30").Value = '  +17.73%  '

$ws.Range("D31").Value = '9.42'
$ws.Range("E31").Value = '  +5.02%  '

$ws.Range("D32").Value = '31.51'
$ws.Range("E32").Value = '  +1.53%  '

$ws.Range("D33").Value = '676.97'
$ws.Range("E33").Value = '  +9.53%  '

$ws.Range("D34").Value = '12.18'
$ws.Range("E34").Value = '  +0.40%  '

$ws.Range("D35").Value = '63.87'
$ws.Range("E35").Value = '  +0.61%  '

$ws.Range("E36").Value = '  +0.54%  '

$ws.Range("D37").Value = '42.43'
$ws.Range("E37").Value = '  +2.53%  '

$ws.Range("D38").Value = '0.429'
$ws.Range("E38").Value = '  +8.95%  '

$ws.Range("D39").Value = '0.999'
$ws.Range("E39").Value = '  -0.19%  '

$ws.Range("D40").Value = '0.0₃0774'
$ws.Range("E40").Value = '  +2.30%  '

$ws.Range("D41").Value = '3.22'
$ws.Range("E41").Value = '  +14.32%  '

$ws.Range("D42").Value = '3.240.09'
$ws.Range("E42").Value = '  +7.35%  '

$ws.Range("E43").Value = '  +3.45%  '

$ws.Range("E44").Value = '  +11.04%  '

$ws.Range("B45").Value = 'dogwifhat'
$ws.Range("C45").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D45").Value = '3.04'
$ws.Range("E45").Value = '  +29.80%  '

$ws.Range("B46").Value = 'FirstDigitalUSD'
$ws.Range("C46").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D46").Value = '0.998'
$ws.Range("E46").Value = '  -0.06%  '

$ws.Range("D47").Value = '0.0417'
$ws.Range("E47").Value = '  +2.35%  '

$ws.Range("E48").Value = '  +10.71%  '

$ws.Range("B49").Value = 'THORChain'
$ws.Range("C49").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D49").Value = '8.82'
$ws.Range("E49").Value = '  +2.56%  '

$ws.Range("D50").Value = '0.131'
$ws.Range("E50").Value = '  +0.44%  '

$ws.Range("B51").Value = 'ApeXProtocol'
$ws.Range("C51").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D51").Value = '3.09'
$ws.Range("E51").Value = '  +1.44%  '
